$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.04837263193803665
$ws.Range("D2").Value = 0.08076544123147755
$ws.Range("E2").Value = 0.1133092308110797
$ws.Range("F2").Value = 3.954288008995718
$ws.Range("G2").Value = 3.517610816729587
$ws.Range("H2").Value = 2.358650277960692
$ws.Range("J2").Value = 0.2323591097577946
$ws.Range("C3").Value = 0.04291586150974069
$ws.Range("D3").Value = 0.07849806491754663
$ws.Range("E3").Value = 0.1129899566129957
$ws.Range("F3").Value = 3.771370646392
$ws.Range("G3").Value = 3.301587764275325
$ws.Range("H3").Value = 2.26737350247339
$ws.Range("J3").Value = 0.2254799030836381
$ws.Range("C4").Value = 0.03959561195027561
$ws.Range("D4").Value = 0.0770883509388085
$ws.Range("E4").Value = 0.1128551022633211
$ws.Range("F4").Value = 3.66221381334239
$ws.Range("G4").Value = 3.171817906806723
$ws.Range("H4").Value = 2.213175233925881
$ws.Range("J4").Value = 0.2214329583296575
$ws.Range("C5").Value = 0.03824977303878541
$ws.Range("D5").Value = 0.07650930721145244
$ws.Range("E5").Value = 0.1128154488327269
$ws.Range("F5").Value = 3.618508684877384
$ws.Range("G5").Value = 3.119638749499074
$ws.Range("H5").Value = 2.191544775453224
$ws.Range("J5").Value = 0.2198275908243659
$ws.Range("C6").Value = 0.0380267198869717
$ws.Range("D6").Value = 0.07641287644713657
$ws.Range("E6").Value = 0.1128097861513808
$ws.Range("F6").Value = 3.61129795376678
$ws.Range("G6").Value = 3.111016377872431
$ws.Range("H6").Value = 2.187980339722913
$ws.Range("J6").Value = 0.2195636460188126
$ws.Range("C7").Value = 0.03957743292282601
$ws.Range("D7").Value = 0.07708056048697642
$ws.Range("E7").Value = 0.1128545056389143
$ws.Range("F7").Value = 3.661621265119948
$ws.Range("G7").Value = 3.17111137751283
$ws.Range("H7").Value = 2.212881683554087
$ws.Range("J7").Value = 0.2214111313567173
$ws.Range("C8").Value = 0.0464846113614783
$ws.Range("D8").Value = 0.07998718910447167
$ws.Range("E8").Value = 0.1131863952433925
$ws.Range("F8").Value = 3.890553302341175
$ws.Range("G8").Value = 3.442519322136206
$ws.Range("H8").Value = 2.326789733958549
$ws.Range("J8").Value = 0.2299500163076971
$ws.Range("C9").Value = 0.06028899607868254
$ws.Range("D9").Value = 0.08555566712294649
$ws.Range("E9").Value = 0.1143267567237061
$ws.Range("F9").Value = 4.365297793557744
$ws.Range("G9").Value = 3.998390302803273
$ws.Range("H9").Value = 2.565213172815845
$ws.Range("J9").Value = 0.2481315412706664
$ws.Range("C10").Value = 0.0706178245743132
$ws.Range("D10").Value = 0.08957788668904954
$ws.Range("E10").Value = 0.1154691160534647
$ws.Range("F10").Value = 4.730963310585878
$ws.Range("G10").Value = 4.42248299043041
$ws.Range("H10").Value = 2.750149006206811
$ws.Range("J10").Value = 0.2624139554791896
$ws.Range("C11").Value = 0.07536367885160189
$ws.Range("D11").Value = 0.0913952357905643
$ws.Range("E11").Value = 0.1160563421569378
$ws.Range("F11").Value = 4.901235078545824
$ws.Range("G11").Value = 4.61911193280207
$ws.Range("H11").Value = 2.836535916717537
$ws.Range("J11").Value = 0.2691231557129612
$ws.Range("C12").Value = 0.07716812685853824
$ws.Range("D12").Value = 0.09208184401928321
$ws.Range("E12").Value = 0.1162885415913735
$ws.Range("F12").Value = 4.966298633670078
$ws.Range("G12").Value = 4.694127045197376
$ws.Range("H12").Value = 2.869584196910068
$ws.Range("J12").Value = 0.271695155198131
$ws.Range("C13").Value = 0.07677917397731449
$ws.Range("D13").Value = 0.0919340378228668
$ws.Range("E13").Value = 0.1162380942381276
$ws.Range("F13").Value = 4.952259668129614
$ws.Range("G13").Value = 4.677946097155598
$ws.Range("H13").Value = 2.86245156676847
$ws.Range("J13").Value = 0.2711398195665708
$ws.Range("C14").Value = 0.07551198269690929
$ws.Range("D14").Value = 0.09145175402433381
$ws.Range("E14").Value = 0.1160752476702456
$ws.Range("F14").Value = 4.906576046229361
$ws.Range("G14").Value = 4.625272194154547
$ws.Range("H14").Value = 2.839248033856109
$ws.Range("J14").Value = 0.2693341218034391
$ws.Range("C15").Value = 0.0747367568904167
$ws.Range("D15").Value = 0.09115614115565052
$ws.Range("E15").Value = 0.1159767830357517
$ws.Range("F15").Value = 4.878670353067548
$ws.Range("G15").Value = 4.59308100855128
$ws.Range("H15").Value = 2.825079204977328
$ws.Range("J15").Value = 0.2682321917641559
$ws.Range("C16").Value = 0.07030866489482435
$ws.Range("D16").Value = 0.08945888632182175
$ws.Range("E16").Value = 0.1154321086230929
$ws.Range("F16").Value = 4.71991643304068
$ws.Range("G16").Value = 4.409709365075059
$ws.Range("H16").Value = 2.744549769859191
$ws.Range("J16").Value = 0.2619798349273026
$ws.Range("C17").Value = 0.06760462282073831
$ws.Range("D17").Value = 0.08841464377898234
$ws.Range("E17").Value = 0.115115352460144
$ws.Range("F17").Value = 4.623546572525697
$ws.Range("G17").Value = 4.298182346519582
$ws.Range("H17").Value = 2.695733525950004
$ws.Range("J17").Value = 0.258199144506392
$ws.Range("C18").Value = 0.06605374256325547
$ws.Range("D18").Value = 0.0878128457814995
$ws.Range("E18").Value = 0.1149395147331163
$ws.Range("F18").Value = 4.568485704440889
$ws.Range("G18").Value = 4.234382464436408
$ws.Range("H18").Value = 2.667867652809377
$ws.Range("J18").Value = 0.2560444970522013
$ws.Range("C19").Value = 0.06552938396440311
$ws.Range("D19").Value = 0.08760887918712257
$ws.Range("E19").Value = 0.1148810666594109
$ws.Range("F19").Value = 4.549905725112012
$ws.Range("G19").Value = 4.212839919092232
$ws.Range("H19").Value = 2.658468829434071
$ws.Range("J19").Value = 0.2553183630229512
$ws.Range("C20").Value = 0.06789201234779796
$ws.Range("D20").Value = 0.08852592561383688
$ws.Range("E20").Value = 0.1151484134686847
$ws.Range("F20").Value = 4.633767022025381
$ws.Range("G20").Value = 4.310018458435138
$ws.Range("H20").Value = 2.700908087636037
$ws.Range("J20").Value = 0.2585995390349467
$ws.Range("C21").Value = 0.0758839853724993
$ws.Range("D21").Value = 0.09159345381031159
$ws.Range("E21").Value = 0.1161228119448339
$ws.Range("F21").Value = 4.91997838171585
$ws.Range("G21").Value = 4.640728520165851
$ws.Range("H21").Value = 2.846054283815135
$ws.Range("J21").Value = 0.2698636404683015
$ws.Range("C22").Value = 0.08114997452028661
$ws.Range("D22").Value = 0.09358913304130567
$ws.Range("E22").Value = 0.1168169890189468
$ws.Range("F22").Value = 5.110455770377826
$ws.Range("G22").Value = 4.860118637852338
$ws.Range("H22").Value = 2.942876066604299
$ws.Range("J22").Value = 0.2774086648392853
$ws.Range("C23").Value = 0.07833533198036946
$ws.Range("D23").Value = 0.09252476986625879
$ws.Range("E23").Value = 0.1164412069167824
$ws.Range("F23").Value = 5.008474398812893
$ws.Range("G23").Value = 4.742720673453732
$ws.Range("H23").Value = 2.891017476649836
$ws.Range("J23").Value = 0.2733646751592573
$ws.Range("C24").Value = 0.06776207202621265
$ws.Range("D24").Value = 0.08847561962447514
$ws.Range("E24").Value = 0.1151334470669205
$ws.Range("F24").Value = 4.629145288575558
$ws.Range("G24").Value = 4.304666361948989
$ws.Range("H24").Value = 2.698568048523214
$ws.Range("J24").Value = 0.2584184617972767
$ws.Range("C25").Value = 0.05652386341895976
$ws.Range("D25").Value = 0.08406217762201607
$ws.Range("E25").Value = 0.1139652353035956
$ws.Range("F25").Value = 4.2339883756853
$ws.Range("G25").Value = 3.845351347943108
$ws.Range("H25").Value = 2.499042109264337
$ws.Range("J25").Value = 0.243054214178386
